# Updated cryptos list on Mon Feb 12 17:24:28 UTC 2024 with GitHub Actions
#
# Applies the per-row Price (column D) / Volume(1h) (column E) updates, plus
# the two row swaps (Toncoin<->Kaspa at rows 28/29, Stellar<->Monero at rows
# 40/41), to the crypto price table on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the (unstyled) default look of the data cells so that after we
# temporarily force a Text number format (to stop Excel from silently
# re-interpreting price strings such as "108.70" or "0.530" as numbers and
# dropping significant trailing zeros) we can restore the cell to its
# original, unstyled appearance.
$normalStyle = $ws.Range("D2").Style

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $normalStyle
}

# Row 2 - Bitcoin
Set-TextValue "D2" "49.888.99"
$ws.Range("E2").Value = "  +3.28%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.561.16"
$ws.Range("E3").Value = "  +1.87%  "

# Row 4 - TetherUSD (price unchanged)
$ws.Range("E4").Value = "  -0.18%  "

# Row 5 - BNB
Set-TextValue "D5" "323.18"
$ws.Range("E5").Value = "  +0.45%  "

# Row 6 - Solana
Set-TextValue "D6" "108.70"
$ws.Range("E6").Value = "  -0.13%  "

# Row 7 - XRP
Set-TextValue "D7" "0.530"
$ws.Range("E7").Value = "  +0.19%  "

# Row 8 - USDC
Set-TextValue "D8" "0.998"
$ws.Range("E8").Value = "  -0.16%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.557"
$ws.Range("E9").Value = "  +2.20%  "

# Row 10 - Avalanche
Set-TextValue "D10" "40.44"
$ws.Range("E10").Value = "  +1.00%  "

# Row 11 - Chainlink
Set-TextValue "D11" "20.41"
$ws.Range("E11").Value = "  +1.50%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.0818"
$ws.Range("E12").Value = "  -0.19%  "

# Row 13 - TRON (price unchanged)
$ws.Range("E13").Value = "  +0.45%  "

# Row 14 - Polkadot (price unchanged)
$ws.Range("E14").Value = "  +0.87%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.961.12"
$ws.Range("E15").Value = "  +1.86%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "2.551.17"
$ws.Range("E16").Value = "  +1.42%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.864"
$ws.Range("E17").Value = "  +2.04%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "49.668.88"
$ws.Range("E18").Value = "  +3.15%  "

# Row 19 - ImmutableX
Set-TextValue "D19" "3.06"
$ws.Range("E19").Value = "  +11.38%  "

# Row 20 - InternetComputer(DFINITY)
Set-TextValue "D20" "13.31"
$ws.Range("E20").Value = "  +1.37%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.70"
$ws.Range("E21").Value = "  +0.07%  "

# Row 22 - ShibaInu
Set-TextValue "D22" "0.0₃0949"
$ws.Range("E22").Value = "  -0.61%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "284.13"
$ws.Range("E23").Value = "  +2.32%  "

# Row 24 - Litecoin
Set-TextValue "D24" "72.46"
$ws.Range("E24").Value = "  +0.29%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "2.54"
$ws.Range("E25").Value = "  -1.05%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "26.42"
$ws.Range("E26").Value = "  +1.92%  "

# Row 27 - Dai (price unchanged)
$ws.Range("E27").Value = "  -0.15%  "

# Row 28/29 - swap: Toncoin <-> Kaspa
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D28" "0.145"
$ws.Range("E28").Value = "  +4.91%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D29" "2.23"
$ws.Range("E29").Value = "  -3.86%  "

# Row 30 - Cosmos
Set-TextValue "D30" "9.89"
$ws.Range("E30").Value = "  +0.59%  "

# Row 31 - InjectiveProtocol
Set-TextValue "D31" "35.47"
$ws.Range("E31").Value = "  +0.08%  "

# Row 32 - OKB
Set-TextValue "D32" "49.46"
$ws.Range("E32").Value = "  +0.49%  "

# Row 33 - Celestia
Set-TextValue "D33" "19.77"
$ws.Range("E33").Value = "  +1.44%  "

# Row 34 - Filecoin
Set-TextValue "D34" "5.41"
$ws.Range("E34").Value = "  +0.80%  "

# Row 35 - FirstDigitalUSD (price unchanged)
$ws.Range("E35").Value = "  -0.31%  "

# Row 36 - Hedera
Set-TextValue "D36" "0.0788"
$ws.Range("E36").Value = "  +0.26%  "

# Row 37 - ARBITRUM
Set-TextValue "D37" "2.03"
$ws.Range("E37").Value = "  +3.53%  "

# Row 38 - RenderToken
Set-TextValue "D38" "4.73"
$ws.Range("E38").Value = "  +1.75%  "

# Row 39 - LidoDAOToken
Set-TextValue "D39" "3.03"
$ws.Range("E39").Value = "  +2.22%  "

# Row 40/41 - swap: Stellar <-> Monero
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D40" "123.05"
$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D41" "0.112"
$ws.Range("E41").Value = "  +0.35%  "

# Row 42 - EnergySwap
Set-TextValue "D42" "22.34"
$ws.Range("E42").Value = "  +3.77%  "

# Row 43 - WEMIXToken (price unchanged)
$ws.Range("E43").Value = "  -0.24%  "

# Row 44 - VeChain (price unchanged)
$ws.Range("E44").Value = "  +2.07%  "

# Row 45 - NEARProtocol (price unchanged)
$ws.Range("E45").Value = "  +5.39%  "

# Row 46 - Maker
Set-TextValue "D46" "2.022.28"
$ws.Range("E46").Value = "  +0.98%  "

# Row 47 - Stacks
Set-TextValue "D47" "2.04"
$ws.Range("E47").Value = "  +9.55%  "

# Row 48 - ApeXProtocol
Set-TextValue "D48" "2.15"
$ws.Range("E48").Value = "  +8.29%  "

# Row 49 - FraxShare
Set-TextValue "D49" "9.04"
$ws.Range("E49").Value = "  +0.14%  "

# Row 50 - THORChain
Set-TextValue "D50" "5.35"
$ws.Range("E50").Value = "  +2.20%  "

# Row 51 - BitcoinSV
Set-TextValue "D51" "81.62"
$ws.Range("E51").Value = "  +1.57%  "
